$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force text format on B2:E51 so numeric-looking strings (e.g. "1.00",
# "0.0000107") are stored as text rather than being auto-converted to numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "43.605.70", "  -0.93%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.298.19", "  +1.93%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.07%  "),
    @(5, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "94.97", "  +9.35%  "),
    @(6, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "268.10", "  -0.93%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.626", "  +1.19%  "),
    @(8, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  -0.10%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.621", "  +2.26%  "),
    @(10, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "44.68", "  -0.92%  "),
    @(11, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0934", "  +0.76%  "),
    @(12, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "8.09", "  +6.52%  "),
    @(13, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.105", "  +0.01%  "),
    @(14, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.646.68", "  +1.77%  "),
    @(15, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "15.26", "  +2.22%  "),
    @(16, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.853", "  +6.69%  "),
    @(17, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.308.85", "  +1.81%  "),
    @(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "43.558.52", "  -0.96%  "),
    @(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000107", "  +3.02%  "),
    @(20, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.30", "  +4.65%  "),
    @(21, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "71.16", "  +0.87%  "),
    @(22, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.27", "  -4.61%  "),
    @(23, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "237.17", "  +1.32%  "),
    @(24, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "9.53", "  +7.72%  "),
    @(25, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.00%  "),
    @(26, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.30", "  +3.93%  "),
    @(27, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.48", "  -1.11%  "),
    @(28, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.31", "  +0.12%  "),
    @(29, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "3.39", "  -4.79%  "),
    @(30, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "38.53", "  -3.02%  "),
    @(31, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "22.33", "  +7.03%  "),
    @(32, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "171.63", "  -1.59%  "),
    @(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0897", "  -0.19%  "),
    @(34, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.53", "  +2.85%  "),
    @(35, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.126", "  +1.64%  "),
    @(36, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0357", "  +2.13%  "),
    @(37, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.108", "  -3.30%  "),
    @(38, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "4.44", "  +1.71%  "),
    @(39, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.46", "  +0.47%  "),
    @(40, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.232", "  +14.36%  "),
    @(41, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.30", "  +4.47%  "),
    @(42, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.34", "  +17.50%  "),
    @(43, "Celestia", "https://coinranking.com/coin/YQcD0lBl7+celestia-tia", "12.05", "  -4.18%  "),
    @(44, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "5.43", "  -0.43%  "),
    @(45, "MultiversX", "https://coinranking.com/coin/omwkOTglq+multiversx-egld", "61.62", "  -3.54%  "),
    @(46, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "9.01", "  +5.28%  "),
    @(47, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.102", "  +3.57%  "),
    @(48, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "100.14", "  -0.61%  "),
    @(49, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.20", "  -0.39%  "),
    @(50, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "2.522.30", "  +1.60%  "),
    @(51, "WOONetwork", "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo", "0.423", "  -1.86%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Restore the original (default/"Normal") cell style now that the text values are set,
# so no lingering explicit text-number-format is left on the cells.
$ws.Range("B2:E51").Style = "Normal"
